$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "nesmadrahim@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, servinaz@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg"
$ws.Range("G3").Value = "mennatulla.medhat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, servinaz@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg"
$ws.Range("G4").Value = "gehanadel@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, manar.montaser@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg"
$ws.Range("G5").Value = "Nadagouda@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg"
$ws.Range("G6").Value = "AbeerRagheb@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg"
$ws.Range("G7").Value = "aya.saeed@med.asu.edu.eg, Arwa.elnagar@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg"
$ws.Range("G8").Value = "amira.m.ibrahim@med.asu.edu.eg, dina.adel@med.asu.edu.eg, Dalia.m.abdelsalam@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg"
$ws.Range("G9").Value = "merna.said@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg"
$ws.Range("G11").Value = "heba.e.mohamed@med.asu.edu.eg, mona.I.hussein@med.asu.edu.eg, mariam.e.mohammed@med.asu.edu.eg"
$ws.Range("G12").Value = "heba.e.mohamed@med.asu.edu.eg, mona.I.hussein@med.asu.edu.eg, mariam.e.mohammed@med.asu.edu.eg"
$ws.Range("G13").Value = "Aya.alaa-eldein@med.asu.edu.eg, marian.samir@med.asu.edu.eg, ManaratAleslam-Zefan@med.asu.edu.eg"
$ws.Range("G15").Value = "Amr-Saeed@med.asu.edu.eg, Walaa.h.ghanima@med.asu.edu.eg"
$ws.Range("G16").Value = "aya.hanafy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, eman.samir@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg"
$ws.Range("G17").Value = "wafaa.ebida@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, marinasorial@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, eman.samir@med.asu.edu.eg"
$ws.Range("G18").Value = "nesmadrahim@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, servinaz@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg"
$ws.Range("G19").Value = "Mohammedeltanany@med.asu.edu.eg, servinaz@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"
$ws.Range("G20").Value = "heba@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, servinaz@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg"
$ws.Range("G21").Value = "Nadagouda@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg"
$ws.Range("G22").Value = "AbeerRagheb@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg"
$ws.Range("G23").Value = "aya.saeed@med.asu.edu.eg, Arwa.elnagar@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg"
$ws.Range("G24").Value = "amira.m.ibrahim@med.asu.edu.eg, dina.adel@med.asu.edu.eg, Dalia.m.abdelsalam@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg"
$ws.Range("G25").Value = "merna.said@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg"
$ws.Range("G27").Value = "heba.e.mohamed@med.asu.edu.eg, mona.I.hussein@med.asu.edu.eg, mariam.e.mohammed@med.asu.edu.eg"
$ws.Range("G28").Value = "heba.e.mohamed@med.asu.edu.eg, mona.I.hussein@med.asu.edu.eg, mariam.e.mohammed@med.asu.edu.eg"
$ws.Range("G29").Value = "Aya.alaa-eldein@med.asu.edu.eg, marian.samir@med.asu.edu.eg, ManaratAleslam-Zefan@med.asu.edu.eg"
$ws.Range("G31").Value = "Amr-Saeed@med.asu.edu.eg, Walaa.h.ghanima@med.asu.edu.eg"
$ws.Range("G32").Value = "aya.hanafy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, eman.samir@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg"
$ws.Range("G33").Value = "wafaa.ebida@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, marinasorial@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, eman.samir@med.asu.edu.eg"
$ws.Range("G34").Value = "nahla.nagiub@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, servinaz@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, gehanadel@med.asu.edu.eg"
$ws.Range("G35").Value = "heba@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, servinaz@med.asu.edu.eg, System, rana.abozaid@med.asu.edu.eg, gehanadel@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"
$ws.Range("G36").Value = "alshimaa.atef@med.asu.edu.egm, heba@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, servinaz@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg"
$ws.Range("G37").Value = "AbeerRagheb@med.asu.edu.eg, Nadagouda@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg"
$ws.Range("G40").Value = "sara_atawia@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, mernamahrous@med.asu.edu.eg, merna.said@med.asu.edu.eg, nahed.mosaad@med.asu.edu.eg, maimustafa@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, nourhan.osama@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg"
$ws.Range("G41").Value = "esraa.mostafa@med.asu.edu.eg, mernamahrous@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, nadia.mostafa@med.asu.edu.eg, dina.adel@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, nourhan.osama@med.asu.edu.eg, Sarah.Abdelmohsen@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg"
$ws.Range("G43").Value = "mona.I.hussein@med.asu.edu.eg, Sarah.Mahdy@med.asu.edu.eg, youstina.gamil@med.asu.edu.eg"
$ws.Range("G44").Value = "mona.I.hussein@med.asu.edu.eg, Sarah.Mahdy@med.asu.edu.eg, youstina.gamil@med.asu.edu.eg"
$ws.Range("G48").Value = "maryam.ashraf@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg"
$ws.Range("G49").Value = "Monica.Eshak@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, eman.samir@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg"
$ws.Range("G50").Value = "nahla.nagiub@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, servinaz@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, gehanadel@med.asu.edu.eg"
$ws.Range("G51").Value = "heba@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, servinaz@med.asu.edu.eg, System, rana.abozaid@med.asu.edu.eg, gehanadel@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"
$ws.Range("G52").Value = "alshimaa.atef@med.asu.edu.egm, heba@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, servinaz@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg"
$ws.Range("G53").Value = "AbeerRagheb@med.asu.edu.eg, Nadagouda@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg"
$ws.Range("G56").Value = "sara_atawia@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, mernamahrous@med.asu.edu.eg, merna.said@med.asu.edu.eg, nahed.mosaad@med.asu.edu.eg, maimustafa@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, nourhan.osama@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg"
$ws.Range("G57").Value = "esraa.mostafa@med.asu.edu.eg, mernamahrous@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, nadia.mostafa@med.asu.edu.eg, dina.adel@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, nourhan.osama@med.asu.edu.eg, Sarah.Abdelmohsen@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg"
$ws.Range("G59").Value = "mona.I.hussein@med.asu.edu.eg, Sarah.Mahdy@med.asu.edu.eg, youstina.gamil@med.asu.edu.eg"
$ws.Range("G60").Value = "mona.I.hussein@med.asu.edu.eg, Sarah.Mahdy@med.asu.edu.eg, youstina.gamil@med.asu.edu.eg"
$ws.Range("G64").Value = "maryam.ashraf@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg"
$ws.Range("G65").Value = "Monica.Eshak@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, eman.samir@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg"
$ws.Range("G66").Value = "nahla.nagiub@med.asu.edu.eg, gehanadel@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, servinaz@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"
$ws.Range("G67").Value = "mennatulla.medhat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, servinaz@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg"
$ws.Range("G68").Value = "alshimaa.atef@med.asu.edu.egm, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"
$ws.Range("G69").Value = "AbeerRagheb@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg"
$ws.Range("G70").Value = "Amera.a.saad@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, Nadagouda@med.asu.edu.eg"
$ws.Range("G72").Value = "sara_atawia@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, mernamahrous@med.asu.edu.eg, merna.said@med.asu.edu.eg, nahed.mosaad@med.asu.edu.eg, maimustafa@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, nourhan.osama@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg"
$ws.Range("G73").Value = "esraa.mostafa@med.asu.edu.eg, merna.said@med.asu.edu.eg, nahed.mosaad@med.asu.edu.eg, Dalia.m.abdelsalam@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg"
$ws.Range("G79").Value = "Amr-Saeed@med.asu.edu.eg, Walaa.h.ghanima@med.asu.edu.eg"
$ws.Range("G80").Value = "aya.hanafy@med.asu.edu.eg, marina_atef@med.asu.edu.eg, Eman_mohamed@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg"
$ws.Range("G81").Value = "wafaa.ebida@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, marinasorial@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, eman.samir@med.asu.edu.eg"
$ws.Range("G82").Value = "nahla.nagiub@med.asu.edu.eg, gehanadel@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, servinaz@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"
$ws.Range("G83").Value = "Mohammedeltanany@med.asu.edu.eg, servinaz@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"
$ws.Range("G84").Value = "heba@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, servinaz@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg"
$ws.Range("G85").Value = "AbeerRagheb@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg"
$ws.Range("G86").Value = "Amera.a.saad@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, Nadagouda@med.asu.edu.eg"
$ws.Range("G88").Value = "sara_atawia@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, mernamahrous@med.asu.edu.eg, merna.said@med.asu.edu.eg, nahed.mosaad@med.asu.edu.eg, maimustafa@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, nourhan.osama@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg"
$ws.Range("G89").Value = "esraa.mostafa@med.asu.edu.eg, merna.said@med.asu.edu.eg, nahed.mosaad@med.asu.edu.eg, Dalia.m.abdelsalam@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg"
$ws.Range("G95").Value = "Amr-Saeed@med.asu.edu.eg, Walaa.h.ghanima@med.asu.edu.eg"
$ws.Range("G96").Value = "aya.hanafy@med.asu.edu.eg, marina_atef@med.asu.edu.eg, Eman_mohamed@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg"
$ws.Range("G98").Value = "nesmadrahim@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, servinaz@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg"
$ws.Range("G100").Value = "alshimaa.atef@med.asu.edu.egm, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"
$ws.Range("G101").Value = "AbeerRagheb@med.asu.edu.eg, Nadagouda@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg"
$ws.Range("G102").Value = "Amera.a.saad@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, Nadagouda@med.asu.edu.eg"
$ws.Range("G103").Value = "aya.saeed@med.asu.edu.eg, Arwa.elnagar@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg"
$ws.Range("G104").Value = "esraa.mostafa@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, dina.adel@med.asu.edu.eg, nourhan.osama@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg"
$ws.Range("G105").Value = "esraa.mostafa@med.asu.edu.eg, mernamahrous@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, nadia.mostafa@med.asu.edu.eg, dina.adel@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, nourhan.osama@med.asu.edu.eg, Sarah.Abdelmohsen@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg"
$ws.Range("G106").Value = "Rania.a.youssef@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"
$ws.Range("G107").Value = "mona.I.hussein@med.asu.edu.eg, Sarah.Mahdy@med.asu.edu.eg, youstina.gamil@med.asu.edu.eg"
$ws.Range("G108").Value = "mona.I.hussein@med.asu.edu.eg, Sarah.Mahdy@med.asu.edu.eg, youstina.gamil@med.asu.edu.eg"
$ws.Range("G112").Value = "youstina.magdy@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, nahla91@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg"
$ws.Range("G114").Value = "nesmadrahim@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, servinaz@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg"
$ws.Range("G116").Value = "gehanadel@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, manar.montaser@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg"
$ws.Range("G117").Value = "AbeerRagheb@med.asu.edu.eg, Nadagouda@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg"
$ws.Range("G118").Value = "Amera.a.saad@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, Nadagouda@med.asu.edu.eg"
$ws.Range("G119").Value = "aya.saeed@med.asu.edu.eg, Arwa.elnagar@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg"
$ws.Range("G120").Value = "esraa.mostafa@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, dina.adel@med.asu.edu.eg, nourhan.osama@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg"
$ws.Range("G121").Value = "esraa.mostafa@med.asu.edu.eg, mernamahrous@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, nadia.mostafa@med.asu.edu.eg, dina.adel@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, nourhan.osama@med.asu.edu.eg, Sarah.Abdelmohsen@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg"
$ws.Range("G122").Value = "Rania.a.youssef@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"
$ws.Range("G123").Value = "mona.I.hussein@med.asu.edu.eg, Sarah.Mahdy@med.asu.edu.eg, youstina.gamil@med.asu.edu.eg"
$ws.Range("G124").Value = "mona.I.hussein@med.asu.edu.eg, Sarah.Mahdy@med.asu.edu.eg, youstina.gamil@med.asu.edu.eg"
$ws.Range("G125").Value = "nancy.abdelshafy@med.asu.edu.eg, Walaa.h.ghanima@med.asu.edu.eg"
$ws.Range("G128").Value = "youstina.magdy@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, nahla91@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg"
